$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.Value = "'30.312.06"
$cell.ClearFormats()

$cell = $ws.Range('D3')
$cell.Value = "'1.871.70"
$cell.ClearFormats()

$cell = $ws.Range('E3')
$cell.Value = "'  +0.50%  "
$cell.ClearFormats()

$cell = $ws.Range('D4')
$cell.Value = "'1.002"
$cell.ClearFormats()

$cell = $ws.Range('E4')
$cell.Value = "'  +0.07%  "
$cell.ClearFormats()

$cell = $ws.Range('D5')
$cell.Value = "'236.03"
$cell.ClearFormats()

$cell = $ws.Range('E5')
$cell.Value = "'  +0.18%  "
$cell.ClearFormats()

$cell = $ws.Range('E6')
$cell.Value = "'  +0.02%  "
$cell.ClearFormats()

$cell = $ws.Range('D7')
$cell.Value = "'0.4710"
$cell.ClearFormats()

$cell = $ws.Range('E7')
$cell.Value = "'  +0.84%  "
$cell.ClearFormats()

$cell = $ws.Range('D8')
$cell.Value = "'0.2886"
$cell.ClearFormats()

$cell = $ws.Range('E8')
$cell.Value = "'  +1.68%  "
$cell.ClearFormats()

$cell = $ws.Range('D9')
$cell.Value = "'0.06631"
$cell.ClearFormats()

$cell = $ws.Range('E9')
$cell.Value = "'  +1.63%  "
$cell.ClearFormats()

$cell = $ws.Range('D10')
$cell.Value = "'21.66"
$cell.ClearFormats()

$cell = $ws.Range('E10')
$cell.Value = "'  +0.00%  "
$cell.ClearFormats()

$cell = $ws.Range('D11')
$cell.Value = "'0.08047"
$cell.ClearFormats()

$cell = $ws.Range('E11')
$cell.Value = "'  +1.36%  "
$cell.ClearFormats()

$cell = $ws.Range('D12')
$cell.Value = "'97.37"
$cell.ClearFormats()

$cell = $ws.Range('E12')
$cell.Value = "'  +0.05%  "
$cell.ClearFormats()

$cell = $ws.Range('D13')
$cell.Value = "'1.873.28"
$cell.ClearFormats()

$cell = $ws.Range('E13')
$cell.Value = "'  +0.47%  "
$cell.ClearFormats()

$cell = $ws.Range('D14')
$cell.Value = "'5.146"
$cell.ClearFormats()

$cell = $ws.Range('E14')
$cell.Value = "'  -0.02%  "
$cell.ClearFormats()

$cell = $ws.Range('D15')
$cell.Value = "'0.6885"
$cell.ClearFormats()

$cell = $ws.Range('E15')
$cell.Value = "'  +1.49%  "
$cell.ClearFormats()

$cell = $ws.Range('D16')
$cell.Value = "'271.85"
$cell.ClearFormats()

$cell = $ws.Range('E16')
$cell.Value = "'  -2.83%  "
$cell.ClearFormats()

$cell = $ws.Range('D17')
$cell.Value = "'30.313.48"
$cell.ClearFormats()

$cell = $ws.Range('E17')
$cell.Value = "'  +0.19%  "
$cell.ClearFormats()

$cell = $ws.Range('D18')
$cell.Value = "'14.20"
$cell.ClearFormats()

$cell = $ws.Range('E18')
$cell.Value = "'  +5.58%  "
$cell.ClearFormats()

$cell = $ws.Range('D19')
$cell.Value = "'0.000007795"
$cell.ClearFormats()

$cell = $ws.Range('E19')
$cell.Value = "'  +6.71%  "
$cell.ClearFormats()

$cell = $ws.Range('E20')
$cell.Value = "'  +0.06%  "
$cell.ClearFormats()

$cell = $ws.Range('D21')
$cell.Value = "'2.117.75"
$cell.ClearFormats()

$cell = $ws.Range('E21')
$cell.Value = "'  +0.18%  "
$cell.ClearFormats()

$cell = $ws.Range('D22')
$cell.Value = "'5.322"
$cell.ClearFormats()

$cell = $ws.Range('E22')
$cell.Value = "'  -1.19%  "
$cell.ClearFormats()

$cell = $ws.Range('D23')
$cell.Value = "'1.001"
$cell.ClearFormats()

$cell = $ws.Range('E23')
$cell.Value = "'  -0.06%  "
$cell.ClearFormats()

$cell = $ws.Range('D24')
$cell.Value = "'6.220"
$cell.ClearFormats()

$cell = $ws.Range('E24')
$cell.Value = "'  +0.79%  "
$cell.ClearFormats()

$cell = $ws.Range('D25')
$cell.Value = "'9.380"
$cell.ClearFormats()

$cell = $ws.Range('E25')
$cell.Value = "'  +2.21%  "
$cell.ClearFormats()

$cell = $ws.Range('D26')
$cell.Value = "'168.16"
$cell.ClearFormats()

$cell = $ws.Range('E26')
$cell.Value = "'  +0.52%  "
$cell.ClearFormats()

$cell = $ws.Range('D27')
$cell.Value = "'18.98"
$cell.ClearFormats()

$cell = $ws.Range('E27')
$cell.Value = "'  -0.48%  "
$cell.ClearFormats()

$cell = $ws.Range('D28')
$cell.Value = "'1.959"
$cell.ClearFormats()

$cell = $ws.Range('E28')
$cell.Value = "'  +1.58%  "
$cell.ClearFormats()

$cell = $ws.Range('E29')
$cell.Value = "'  -1.00%  "
$cell.ClearFormats()

$cell = $ws.Range('D30')
$cell.Value = "'0.09959"
$cell.ClearFormats()

$cell = $ws.Range('E30')
$cell.Value = "'  +2.33%  "
$cell.ClearFormats()

$cell = $ws.Range('D31')
$cell.Value = "'4.371"
$cell.ClearFormats()

$cell = $ws.Range('E31')
$cell.Value = "'  +0.07%  "
$cell.ClearFormats()

$cell = $ws.Range('E32')
$cell.Value = "'  -0.76%  "
$cell.ClearFormats()

$cell = $ws.Range('D33')
$cell.Value = "'4.089"
$cell.ClearFormats()

$cell = $ws.Range('D34')
$cell.Value = "'0.04711"
$cell.ClearFormats()

$cell = $ws.Range('E34')
$cell.Value = "'  -0.59%  "
$cell.ClearFormats()

$cell = $ws.Range('D35')
$cell.Value = "'1.138"
$cell.ClearFormats()

$cell = $ws.Range('E35')
$cell.Value = "'  +0.74%  "
$cell.ClearFormats()

$cell = $ws.Range('D36')
$cell.Value = "'0.7021"
$cell.ClearFormats()

$cell = $ws.Range('E36')
$cell.Value = "'  -0.53%  "
$cell.ClearFormats()

$cell = $ws.Range('D37')
$cell.Value = "'2.713"
$cell.ClearFormats()

$cell = $ws.Range('E37')
$cell.Value = "'  +0.09%  "
$cell.ClearFormats()

$cell = $ws.Range('D38')
$cell.Value = "'0.01889"
$cell.ClearFormats()

$cell = $ws.Range('E38')
$cell.Value = "'  +1.19%  "
$cell.ClearFormats()

$cell = $ws.Range('D39')
$cell.Value = "'2.652"
$cell.ClearFormats()

$cell = $ws.Range('E39')
$cell.Value = "'  +2.65%  "
$cell.ClearFormats()

$cell = $ws.Range('D40')
$cell.Value = "'6.314"
$cell.ClearFormats()

$cell = $ws.Range('E40')
$cell.Value = "'  +0.19%  "
$cell.ClearFormats()

$cell = $ws.Range('D41')
$cell.Value = "'72.72"
$cell.ClearFormats()

$cell = $ws.Range('E41')
$cell.Value = "'  -3.06%  "
$cell.ClearFormats()

$cell = $ws.Range('D42')
$cell.Value = "'1.964"
$cell.ClearFormats()

$cell = $ws.Range('E42')
$cell.Value = "'  +0.28%  "
$cell.ClearFormats()

$cell = $ws.Range('D43')
$cell.Value = "'0.8438"
$cell.ClearFormats()

$cell = $ws.Range('E43')
$cell.Value = "'  -0.82%  "
$cell.ClearFormats()

$cell = $ws.Range('D44')
$cell.Value = "'0.4177"
$cell.ClearFormats()

$cell = $ws.Range('E44')
$cell.Value = "'  -0.10%  "
$cell.ClearFormats()

$cell = $ws.Range('D45')
$cell.Value = "'1.000"
$cell.ClearFormats()

$cell = $ws.Range('E45')
$cell.Value = "'  +0.00%  "
$cell.ClearFormats()

$cell = $ws.Range('D46')
$cell.Value = "'103.15"
$cell.ClearFormats()

$cell = $ws.Range('E46')
$cell.Value = "'  -0.33%  "
$cell.ClearFormats()

$cell = $ws.Range('D47')
$cell.Value = "'9.330"
$cell.ClearFormats()

$cell = $ws.Range('E47')
$cell.Value = "'  +0.35%  "
$cell.ClearFormats()

$cell = $ws.Range('D48')
$cell.Value = "'7.100"
$cell.ClearFormats()

$cell = $ws.Range('E48')
$cell.Value = "'  -0.96%  "
$cell.ClearFormats()

$cell = $ws.Range('D49')
$cell.Value = "'933.02"
$cell.ClearFormats()

$cell = $ws.Range('E49')
$cell.Value = "'  -3.27%  "
$cell.ClearFormats()

$cell = $ws.Range('D50')
$cell.Value = "'34.52"
$cell.ClearFormats()

$cell = $ws.Range('E50')
$cell.Value = "'  +1.39%  "
$cell.ClearFormats()

$cell = $ws.Range('D51')
$cell.Value = "'0.05685"
$cell.ClearFormats()

$cell = $ws.Range('E51')
$cell.Value = "'  +0.65%  "
$cell.ClearFormats()
